$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells (coin name, link, volume label) - assign directly
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E14").Value = "13BitForexTokenBF"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E15").Value = "14CoinExTokenCET"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("E16").Value = "15TigerCashTCH"
$ws.Range("B17").Value = "HotbitToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("E17").Value = "16HotbitTokenHTB"
$ws.Range("B18").Value = "BitKan"
$ws.Range("C18").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("E18").Value = "17BitKanKAN"
$ws.Range("B19").Value = "NitroEx"
$ws.Range("C19").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("E19").Value = "18NitroExNTX"
$ws.Range("B20").Value = "LEO"
$ws.Range("C20").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E20").Value = "19LEOLEO"
$ws.Range("B21").Value = "KuCoinToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("E21").Value = "20KuCoinTokenKCS"
$ws.Range("B22").Value = "BTSEToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("E22").Value = "21BTSETokenBTSE"
$ws.Range("B23").Value = "One"
$ws.Range("C23").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E23").Value = "22OneONEBestin24h"
$ws.Range("B24").Value = "BitpandaEcosystemToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("E24").Value = "23BitpandaEcosystemTokenBEST"
$ws.Range("B25").Value = "ProBitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("E25").Value = "24ProBitTokenPROB"
$ws.Range("B26").Value = "MCDex"
$ws.Range("C26").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E26").Value = "25MCDexMCB"
$ws.Range("E43").Value = "42CEJICEJI"

# Numeric-looking text cells (price, hour) - force text storage so the
# literal string (incl. trailing zeros) round-trips, then restore General format
$numericTextCells = @("D2","G2","D3","G3","D4","G4","D5","G5","D6","G6","D7","G7","D8","G8","D9","G9","D10","G10","D11","G11","D12","G12","D13","G13","D14","G14","D15","G15","D16","G16","D17","G17","D18","G18","D19","G19","D20","G20","D21","G21","D22","G22","D23","G23","D24","G24","D25","G25","D26","G26","D27","G27","G28","G29","G30","G31","G32","G33","G34","G35","G36","G37","G38","G39","G40","D41","G41","D42","G42","G43","D44","G44","G45","G46","D47","G47","D48","G48","G49","G50","G51")
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}
$ws.Range("D2").Value = "246.19"
$ws.Range("G2").Value = "9"
$ws.Range("D3").Value = "22.81"
$ws.Range("G3").Value = "9"
$ws.Range("D4").Value = "5.399"
$ws.Range("G4").Value = "9"
$ws.Range("D5").Value = "0.05744"
$ws.Range("G5").Value = "9"
$ws.Range("D6").Value = "3.434"
$ws.Range("G6").Value = "9"
$ws.Range("D7").Value = "0.8138"
$ws.Range("G7").Value = "9"
$ws.Range("D8").Value = "0.8955"
$ws.Range("G8").Value = "9"
$ws.Range("D9").Value = "0.1442"
$ws.Range("G9").Value = "9"
$ws.Range("D10").Value = "0.07371"
$ws.Range("G10").Value = "9"
$ws.Range("D11").Value = "0.02990"
$ws.Range("G11").Value = "9"
$ws.Range("D12").Value = "0.03068"
$ws.Range("G12").Value = "9"
$ws.Range("D13").Value = "0.09408"
$ws.Range("G13").Value = "9"
$ws.Range("D14").Value = "0.001575"
$ws.Range("G14").Value = "9"
$ws.Range("D15").Value = "0.04826"
$ws.Range("G15").Value = "9"
$ws.Range("D16").Value = "0.006159"
$ws.Range("G16").Value = "9"
$ws.Range("D17").Value = "0.005111"
$ws.Range("G17").Value = "9"
$ws.Range("D18").Value = "0.0009950"
$ws.Range("G18").Value = "9"
$ws.Range("D19").Value = "0.0001500"
$ws.Range("G19").Value = "9"
$ws.Range("D20").Value = "3.746"
$ws.Range("G20").Value = "9"
$ws.Range("D21").Value = "6.331"
$ws.Range("G21").Value = "9"
$ws.Range("D22").Value = "2.199"
$ws.Range("G22").Value = "9"
$ws.Range("D23").Value = "0.01092"
$ws.Range("G23").Value = "9"
$ws.Range("D24").Value = "0.3275"
$ws.Range("G24").Value = "9"
$ws.Range("D25").Value = "0.1310"
$ws.Range("G25").Value = "9"
$ws.Range("D26").Value = "4.157"
$ws.Range("G26").Value = "9"
$ws.Range("D27").Value = "0.0003158"
$ws.Range("G27").Value = "9"
$ws.Range("G28").Value = "9"
$ws.Range("G29").Value = "9"
$ws.Range("G30").Value = "9"
$ws.Range("G31").Value = "9"
$ws.Range("G32").Value = "9"
$ws.Range("G33").Value = "9"
$ws.Range("G34").Value = "9"
$ws.Range("G35").Value = "9"
$ws.Range("G36").Value = "9"
$ws.Range("G37").Value = "9"
$ws.Range("G38").Value = "9"
$ws.Range("G39").Value = "9"
$ws.Range("G40").Value = "9"
$ws.Range("D41").Value = "0.006783"
$ws.Range("G41").Value = "9"
$ws.Range("D42").Value = "0.1069"
$ws.Range("G42").Value = "9"
$ws.Range("G43").Value = "9"
$ws.Range("D44").Value = "0.007373"
$ws.Range("G44").Value = "9"
$ws.Range("G45").Value = "9"
$ws.Range("G46").Value = "9"
$ws.Range("D47").Value = "0.3799"
$ws.Range("G47").Value = "9"
$ws.Range("D48").Value = "0.1689"
$ws.Range("G48").Value = "9"
$ws.Range("G49").Value = "9"
$ws.Range("G50").Value = "9"
$ws.Range("G51").Value = "9"
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).NumberFormat = "General"
}
